$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their text formatting (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "0.0000282" into numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.848.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.415.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +7.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.186"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.593"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.87"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.65%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.963.30"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "648.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +12.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.899.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.422.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.919"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.05"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.82"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.80%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "612.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.053.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.16%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.64%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.33"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.347"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0427"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.87"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.30%  "
